$wb = $excel.ActiveWorkbook

# -- sheet2: header row D1:K1 changes from the generic "Name" string to
#    distinct Name1..Name8 strings (adds 8 new shared strings) --
$ws2 = $wb.Worksheets.Item("sheet2")
$ws2.Range("D1").Value = "Name1"
$ws2.Range("E1").Value = "Name2"
$ws2.Range("F1").Value = "Name3"
$ws2.Range("G1").Value = "Name4"
$ws2.Range("H1").Value = "Name5"
$ws2.Range("I1").Value = "Name6"
$ws2.Range("J1").Value = "Name7"
$ws2.Range("K1").Value = "Name8"

# -- make "sheet2" the active tab/sheet (was "Test") and move its
#    selection from C1:C2 to L11 --
$ws2.Activate()
$null = $ws2.Range("L11").Select()
